$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "96-3="
$t.Cell(1,2).Range.Text = "4+46="
$t.Cell(1,3).Range.Text = "27+33="
$t.Cell(1,4).Range.Text = "27+72="
$t.Cell(1,5).Range.Text = "0+67="
$t.Cell(2,1).Range.Text = "77-20="
$t.Cell(2,2).Range.Text = "36-7="
$t.Cell(2,3).Range.Text = "14+68="
$t.Cell(2,4).Range.Text = "6+58="
$t.Cell(2,5).Range.Text = "16+71="
$t.Cell(3,1).Range.Text = "10+28="
$t.Cell(3,2).Range.Text = "89-46="
$t.Cell(3,3).Range.Text = "95-22="
$t.Cell(3,4).Range.Text = "86-42="
$t.Cell(3,5).Range.Text = "32+1="
$t.Cell(4,1).Range.Text = "6+84="
$t.Cell(4,2).Range.Text = "21+10="
$t.Cell(4,3).Range.Text = "99-33="
$t.Cell(4,4).Range.Text = "56+17="
$t.Cell(4,5).Range.Text = "13+73="
$t.Cell(5,1).Range.Text = "36-7="
$t.Cell(5,2).Range.Text = "25+1="
$t.Cell(5,3).Range.Text = "39-22="
$t.Cell(5,4).Range.Text = "52-18="
$t.Cell(5,5).Range.Text = "94-39="
$t.Cell(6,1).Range.Text = "64-26="
$t.Cell(6,2).Range.Text = "6+58="
$t.Cell(6,3).Range.Text = "19+52="
$t.Cell(6,4).Range.Text = "33+39="
$t.Cell(6,5).Range.Text = "62-48="
$t.Cell(7,1).Range.Text = "1+79="
$t.Cell(7,2).Range.Text = "59+22="
$t.Cell(7,3).Range.Text = "91-85="
$t.Cell(7,4).Range.Text = "9+46="
$t.Cell(7,5).Range.Text = "73-40="
$t.Cell(8,1).Range.Text = "67-42="
$t.Cell(8,2).Range.Text = "71+9="
$t.Cell(8,3).Range.Text = "86-2="
$t.Cell(8,4).Range.Text = "70-45="
$t.Cell(8,5).Range.Text = "47+40="
$t.Cell(9,1).Range.Text = "26+42="
$t.Cell(9,2).Range.Text = "52-12="
$t.Cell(9,3).Range.Text = "34-19="
$t.Cell(9,4).Range.Text = "63-4="
$t.Cell(9,5).Range.Text = "56+0="
$t.Cell(10,1).Range.Text = "90-71="
$t.Cell(10,2).Range.Text = "44-26="
$t.Cell(10,3).Range.Text = "70-12="
$t.Cell(10,4).Range.Text = "50-38="
$t.Cell(10,5).Range.Text = "63-51="
$t.Cell(11,1).Range.Text = "87-5="
$t.Cell(11,2).Range.Text = "80-2="
$t.Cell(11,3).Range.Text = "86+5="
$t.Cell(11,4).Range.Text = "34+51="
$t.Cell(11,5).Range.Text = "1+60="
$t.Cell(12,1).Range.Text = "90-44="
$t.Cell(12,2).Range.Text = "72-30="
$t.Cell(12,3).Range.Text = "23+39="
$t.Cell(12,4).Range.Text = "15+18="
$t.Cell(12,5).Range.Text = "67+29="
$t.Cell(13,1).Range.Text = "77+10="
$t.Cell(13,2).Range.Text = "47-40="
$t.Cell(13,3).Range.Text = "84-15="
$t.Cell(13,4).Range.Text = "29-18="
$t.Cell(13,5).Range.Text = "9-5="
$t.Cell(14,1).Range.Text = "5+1="
$t.Cell(14,2).Range.Text = "83-71="
$t.Cell(14,3).Range.Text = "2+50="
$t.Cell(14,4).Range.Text = "89-23="
$t.Cell(14,5).Range.Text = "46-46="
$t.Cell(15,1).Range.Text = "22+0="
$t.Cell(15,2).Range.Text = "80-75="
$t.Cell(15,3).Range.Text = "78-6="
$t.Cell(15,4).Range.Text = "30-21="
$t.Cell(15,5).Range.Text = "72+14="
$t.Cell(16,1).Range.Text = "87+4="
$t.Cell(16,2).Range.Text = "59-1="
$t.Cell(16,3).Range.Text = "72-0="
$t.Cell(16,4).Range.Text = "78+1="
$t.Cell(16,5).Range.Text = "40-33="
$t.Cell(17,1).Range.Text = "14+82="
$t.Cell(17,2).Range.Text = "90-66="
$t.Cell(17,3).Range.Text = "0+18="
$t.Cell(17,4).Range.Text = "67-0="
$t.Cell(17,5).Range.Text = "70-56="
$t.Cell(18,1).Range.Text = "15+15="
$t.Cell(18,2).Range.Text = "98-46="
$t.Cell(18,3).Range.Text = "42-6="
$t.Cell(18,4).Range.Text = "49+38="
$t.Cell(18,5).Range.Text = "97-2="
$t.Cell(19,1).Range.Text = "68+27="
$t.Cell(19,2).Range.Text = "55-39="
$t.Cell(19,3).Range.Text = "64-11="
$t.Cell(19,4).Range.Text = "28-10="
$t.Cell(19,5).Range.Text = "20+1="
$t.Cell(20,1).Range.Text = "2+82="
$t.Cell(20,2).Range.Text = "74+23="
$t.Cell(20,3).Range.Text = "66-32="
$t.Cell(20,4).Range.Text = "50+34="
$t.Cell(20,5).Range.Text = "54+39="
